$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header column, matching the style of the existing header (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new data values for the Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
